$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.951.75"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "3.305.21"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "186.97"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").Value = "583.14"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "3.879.45"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("D13").Value = "0.138"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").Value = "27.50"
$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("D15").Value = "68.105.40"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "3.299.92"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "443.42"
$ws.Range("E18").Value = "  +10.51%  "

$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("D20").Value = "13.58"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("E21").Value = "  +3.06%  "

$ws.Range("D22").Value = "74.51"
$ws.Range("E22").Value = "  +5.12%  "

$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.515"
$ws.Range("E24").Value = "  +1.43%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.453.40"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("E27").Value = "  +1.07%  "

$ws.Range("D28").Value = "9.13"
$ws.Range("E28").Value = "  -4.01%  "

$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("E30").Value = "  +1.87%  "

$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").Value = "1.54"
$ws.Range("E36").Value = "  +6.04%  "

$ws.Range("D37").Value = "163.56"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("D39").Value = "27.17"
$ws.Range("E39").Value = "  +0.73%  "

$ws.Range("D40").Value = "4.52"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("D41").Value = "0.783"
$ws.Range("E41").Value = "  -2.80%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "6.47"
$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.739.40"
$ws.Range("E43").Value = "  +2.65%  "

$ws.Range("D44").Value = "40.47"
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").Value = "0.0675"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.41"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "24.69"
$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("D48").Value = "328.54"
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").Value = "31.89"
$ws.Range("E50").Value = "  +4.42%  "

$ws.Range("D51").Value = "0.994"
$ws.Range("E51").Value = "  +2.69%  "
